# Updates crypto price/volume figures (and the Kaspa / Binance-PegBSC-USD row order)
# to match the latest scrape, per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.627.10"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").Value = "3.707.37"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "673.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +1.08%  "

# Row 9
$ws.Range("E9").Value = "  +0.76%  "

# Row 10
$ws.Range("E10").Value = "  +2.34%  "

# Row 11
$ws.Range("E11").Value = "  +1.86%  "

# Row 12
$ws.Range("E12").Value = "  +1.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("D14").Value = "3.705.10"
$ws.Range("E14").Value = "  +0.65%  "

# Row 15
$ws.Range("D15").Value = "69.671.43"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16
$ws.Range("E16").Value = "  +1.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.96%  "

# Row 18
$ws.Range("E18").Value = "  +2.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "473.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.655"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23
$ws.Range("D23").Value = "3.856.71"
$ws.Range("E23").Value = "  +0.92%  "

# Row 24
$ws.Range("E24").Value = "  +5.67%  "

# Row 25
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("E29").Value = "  +0.33%  "

# Row 30
$ws.Range("E30").Value = "  +1.38%  "

# Row 31
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "

# Row 33
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.167"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "

# Row 35
$ws.Range("D35").Value = "3.695.83"
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "

# Row 39
$ws.Range("E39").Value = "  +2.22%  "

# Row 40
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("E41").Value = "  +1.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
$ws.Range("E44").Value = "  -1.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

# Row 47
$ws.Range("E47").Value = "  +1.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.63%  "

# Row 49
$ws.Range("E49").Value = "  -0.25%  "

# Row 50
$ws.Range("E50").Value = "  +1.81%  "

# Row 51
$ws.Range("E51").Value = "  +1.69%  "
